$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("K2").Value = 1
$ws.Range("N2").Value = 0.1
$ws.Range("O2").Value = 0.1

# Fill in row 3
$ws.Range("B3").Value = 120
$ws.Range("C3").Value = "n"
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 30
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 60
$ws.Range("I3").Value = 90
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 500
$ws.Range("N3").Value = 0.4
$ws.Range("O3").Value = 0.4

# Fill in row 4
$ws.Range("B4").Value = 120
$ws.Range("C4").Value = "n"
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 30
$ws.Range("G4").Value = 60
$ws.Range("H4").Value = 60
$ws.Range("I4").Value = 90
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 500
$ws.Range("N4").Value = 0.7
$ws.Range("O4").Value = 0.7

# Fill in row 5
$ws.Range("B5").Value = 120
$ws.Range("C5").Value = "n"
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 30
$ws.Range("G5").Value = 60
$ws.Range("H5").Value = 60
$ws.Range("I5").Value = 90
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 500
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 1

# Update active cell selection
$ws.Range("D24").Select()
